# Generate Report for Handoff
#
# Updates the localization-status report after a fresh handoff-report run:
#   - "Priority" column (E) is stamped "ht" for the rows that now have a
#     handoff priority, on both the zh-cn and de-de sheets.
#   - The "Latest Handoff/HO Xliff Generate" datetime columns are refreshed
#     to the new report-generation timestamps, on all three sheets.

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 11, 12, 13, 14)

# --- zh-cn sheet -----------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $zhcn.Cells.Item($r, 5).Value = "ht"                       # E: Priority
    $zhcn.Cells.Item($r, 8).Value = "2016-08-23 00:20:52"      # H: Latest Handoff Datetime
}

# --- de-de sheet -------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $dede.Cells.Item($r, 5).Value = "ht"                       # E: Priority
    $dede.Cells.Item($r, 8).Value = "2016-08-23 00:20:57"      # H: Latest Handoff Datetime
}

# --- Overview sheet ------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $overview.Cells.Item($r, 7).Value = "2016-08-23 00:20:57"  # G: Latest HO Xliff Generate Date
}
